# Rename component headers in the "default" worksheet to reflect the
# updated PM2 component names:
#   X_CH  -> X_PG
#   X_LI  -> X_TAG
#   S_F   -> S_G
#
# Also update the current cell selection to J1 (single cell) as captured
# in the workbook when it was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("default")

$ws.Range("D1").Value = "X_PG"
$ws.Range("E1").Value = "X_TAG"
$ws.Range("H1").Value = "S_G"

$ws.Activate()
$ws.Range("J1").Select()
